# arreglo de problemas en gestion productos JAGL
# - corrige el precio en USD (mayor precision) de la fila 39 y agrega
#   una etiqueta vacia en la columna H
# - agrega un nuevo producto en la fila 40

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39: recompute USD price with full precision and tag the row ---
$ws.Range("C39").Value = 179.35897435897436
$ws.Range("D39").Value = 1399
$ws.Range("E39").Value = 9
$ws.Range("H39").Value = ""

# --- Row 40: new product entry ---
$ws.Range("A40").Value = 40
$ws.Range("B40").Value = "afaf"
$ws.Range("C40").Value = 200
$ws.Range("D40").Value = 1560
$ws.Range("E40").Value = 65
$ws.Range("F40").Value = "CONTROL DE ACCESO Y SEGURIDAD"
$ws.Range("G40").Value = "asd"
$ws.Range("H40").Value = "sfsf"
$ws.Range("I40").Value = "🔒"
